$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()

$cell1 = $newRow.Cells.Item(1).Range
$cell1.Text = "15/11/2021"

$cell2 = $newRow.Cells.Item(2).Range
$cell2.Text = "7 Hours"

$cell3 = $newRow.Cells.Item(3).Range
$cell3.Text = "World Generation " + [char]8211 + " Objective 1"

$cell4 = $newRow.Cells.Item(4).Range
$cell4.Text = "Attempted to implement a perlin noise algorithm, to limited success. The algorithm still needs work to be used in the project."
